$wb = $excel.ActiveWorkbook

# ===== Sheet: Log Rating Torneo 02 2022 =====
$ws = $wb.Worksheets.Item("Log Rating Torneo 02 2022")
# swap rows 13 and 14 (cols B:H)
$ws.Range("B13").Value = "Bonelli, Marcos (951)"
$ws.Range("B14").Value = "Migoni, Nicolas (1535)"
$ws.Range("C13").Value = "'601"
$ws.Range("C14").Value = "'17"
$ws.Range("D13").Value = "'1"
$ws.Range("D14").Value = "'18"
$ws.Range("E13").Value = "'-1"
$ws.Range("E14").Value = "'-18"
$ws.Range("H13").Value = "'0.5"
$ws.Range("H14").Value = "'1"

# swap rows 19 and 20 (cols B:H)
$ws.Range("B19").Value = "Nohara, Andres (1375)"
$ws.Range("B20").Value = "Palamedi, Cristian (966)"
$ws.Range("C19").Value = "'602"
$ws.Range("C20").Value = "'1011"
$ws.Range("D19").Value = "'1"
$ws.Range("D20").Value = "'0"
$ws.Range("E19").Value = "'-1"
$ws.Range("E20").Value = "'0"

# swap rows 32 and 33 (cols B:H)
$ws.Range("B32").Value = "Acosta, Gaston (478)"
$ws.Range("B33").Value = "Levin, Raul (994)"
$ws.Range("C32").Value = "'816"
$ws.Range("C33").Value = "'300"
$ws.Range("D32").Value = "'1"
$ws.Range("D33").Value = "'6"
$ws.Range("E32").Value = "'-1"
$ws.Range("E33").Value = "'-6"
$ws.Range("H32").Value = "'0.5"
$ws.Range("H33").Value = "'1"

# swap rows 36 and 37 (cols B:H)
$ws.Range("B36").Value = "Escudero, Martin (956)"
$ws.Range("B37").Value = "Goy, Gerardo (954)"
$ws.Range("C36").Value = "'419"
$ws.Range("C37").Value = "'421"

# swap rows 38 and 39 (cols B:H)
$ws.Range("B38").Value = "Bonelli, Marcos (951)"
$ws.Range("B39").Value = "Palamedi, Cristian (966)"
$ws.Range("C38").Value = "'411"
$ws.Range("C39").Value = "'396"
$ws.Range("D38").Value = "'4"
$ws.Range("D39").Value = "'6"
$ws.Range("E38").Value = "'-4"
$ws.Range("E39").Value = "'-6"

# swap rows 58 and 59 (cols B:H)
$ws.Range("B58").Value = "Larrosa, Jorge (525)"
$ws.Range("B59").Value = "Migoni, Anibal (1010)"
$ws.Range("C58").Value = "'-16"
$ws.Range("C59").Value = "'-501"
$ws.Range("D58").Value = "'18"
$ws.Range("D59").Value = "'64"
$ws.Range("E58").Value = "'-18"
$ws.Range("E59").Value = "'-64"

# ===== Sheet: Partidos Torneo 02 2022 =====
$ws = $wb.Worksheets.Item("Partidos Torneo 02 2022")
# swap rows 20 and 21 (cols A:F)
$ws.Range("B20").Value = "Nohara, Andres"
$ws.Range("B21").Value = "Palamedi, Cristian"

# swap rows 36 and 37 (cols A:F)
$ws.Range("B36").Value = "Escudero, Martin"
$ws.Range("B37").Value = "Goy, Gerardo"

# swap rows 62 and 63 (cols A:F)
$ws.Range("B62").Value = "Benicio, Oscar"
$ws.Range("B63").Value = "Prettis, Juan"

# ===== Sheet: Nivel de Juego Torneo 00 2022 =====
$ws = $wb.Worksheets.Item("Nivel de Juego Torneo 00 2022")
# reorder tied-rating group rows [22, 23] alphabetically by player name
$ws.Range("B22").Value = "Dupertuis, Gaston"
$ws.Range("C22").Value = "Parana"
$ws.Range("D22").Value = "Aspatem"
$ws.Range("B23").Value = "Nowotny, Martin"
$ws.Range("C23").Value = "Libertador San Martin"
$ws.Range("D23").Value = ""

# reorder tied-rating group rows [27, 28] alphabetically by player name
$ws.Range("B27").Value = "Aguirre, German"
$ws.Range("C27").Value = "Parana"
$ws.Range("D27").Value = "Aspatem"
$ws.Range("B28").Value = "La Barba, Pablo"
$ws.Range("C28").Value = "Parana"
$ws.Range("D28").Value = "Aspatem"

# reorder tied-rating group rows [62, 63] alphabetically by player name
$ws.Range("B62").Value = "Sartor, Yemel"
$ws.Range("C62").Value = "Avellaneda"
$ws.Range("D62").Value = "ATMAR"
$ws.Range("B63").Value = "Sueldo, Pablo"
$ws.Range("C63").Value = ""
$ws.Range("D63").Value = ""

# reorder tied-rating group rows [76, 77] alphabetically by player name
$ws.Range("B76").Value = "Becker, Fernando"
$ws.Range("C76").Value = "Parana"
$ws.Range("D76").Value = "Aspatem"
$ws.Range("B77").Value = "Maerker, Shion"
$ws.Range("C77").Value = "Libertador San Martin"
$ws.Range("D77").Value = "CRL"

# reorder tied-rating group rows [80, 81] alphabetically by player name
$ws.Range("B80").Value = "Campos, Dario"
$ws.Range("C80").Value = "Santa Fe"
$ws.Range("D80").Value = "Atemeli"
$ws.Range("B81").Value = "Chiara, Lucio"
$ws.Range("C81").Value = "Parana"
$ws.Range("D81").Value = "Aspatem"

# reorder tied-rating group rows [82, 83] alphabetically by player name
$ws.Range("B82").Value = "Gimenez, Maximo"
$ws.Range("C82").Value = "Parana"
$ws.Range("D82").Value = "Aspatem"
$ws.Range("B83").Value = "Lerch, Juan Carlos"
$ws.Range("C83").Value = "Parana"
$ws.Range("D83").Value = "Tiro Federal"

# reorder tied-rating group rows [102, 103, 104, 105, 106, 107] alphabetically by player name
$ws.Range("B102").Value = "Arrieta, Maximiliano"
$ws.Range("C102").Value = "Libertador San Martin"
$ws.Range("D102").Value = "CRL"
$ws.Range("B103").Value = "Asenie, Santiago"
$ws.Range("C103").Value = "Libertador San Martin"
$ws.Range("D103").Value = ""
$ws.Range("B104").Value = "Badano, Pablo"
$ws.Range("C104").Value = "Parana"
$ws.Range("D104").Value = "Aspatem"
$ws.Range("B106").Value = "Mendieta, Elias"
$ws.Range("C106").Value = "Santa Fe"
$ws.Range("D106").Value = ""
$ws.Range("B107").Value = "Presel, Raul"
$ws.Range("C107").Value = "Parana"
$ws.Range("D107").Value = "Aspatem"

# reorder tied-rating group rows [169, 170] alphabetically by player name
$ws.Range("B169").Value = "Aguirre, Gabriel"
$ws.Range("C169").Value = "Parana"
$ws.Range("D169").Value = "Aspatem"
$ws.Range("B170").Value = "Martinez, Dylan"
$ws.Range("C170").Value = ""
$ws.Range("D170").Value = ""

# reorder tied-rating group rows [175, 176] alphabetically by player name
$ws.Range("B175").Value = "Rodriguez Alarcon Emiliano"
$ws.Range("C175").Value = ""
$ws.Range("D175").Value = ""
$ws.Range("B176").Value = "Werner, Graciela"
$ws.Range("C176").Value = "Parana"
$ws.Range("D176").Value = "Aspatem"

# reorder tied-rating group rows [177, 178, 179] alphabetically by player name
$ws.Range("B177").Value = "Javita, Luis"
$ws.Range("C177").Value = "Libertador San Martin"
$ws.Range("D177").Value = "CRL"
$ws.Range("B179").Value = "Mir, Tomas"
$ws.Range("C179").Value = ""
$ws.Range("D179").Value = ""

# reorder tied-rating group rows [189, 190] alphabetically by player name
$ws.Range("B189").Value = "Jose"
$ws.Range("C189").Value = "Libertador San Martin"
$ws.Range("D189").Value = "CRL"
$ws.Range("B190").Value = "Michea, Ignacio"
$ws.Range("C190").Value = ""
$ws.Range("D190").Value = ""

# reorder tied-rating group rows [194, 195, 196, 197] alphabetically by player name
$ws.Range("B194").Value = "Arrieta, Matias"
$ws.Range("C194").Value = "Libertador San Martin"
$ws.Range("D194").Value = "CRL"
$ws.Range("B196").Value = "Muller, Tomas"
$ws.Range("C196").Value = ""
$ws.Range("D196").Value = ""
$ws.Range("B197").Value = "Portillo, Lucas"
$ws.Range("C197").Value = "Esperanza"
$ws.Range("D197").Value = "ATME"

# reorder tied-rating group rows [202, 203] alphabetically by player name
$ws.Range("B202").Value = "Comas, Javier"
$ws.Range("C202").Value = "Parana"
$ws.Range("D202").Value = "Aspatem"
$ws.Range("B203").Value = "Pillac, Juan Pablo"
$ws.Range("C203").Value = "Libertador San Martin"
$ws.Range("D203").Value = "CRL"

# reorder tied-rating group rows [212, 213] alphabetically by player name
$ws.Range("B212").Value = "Godano, Lucas"
$ws.Range("C212").Value = ""
$ws.Range("D212").Value = ""
$ws.Range("B213").Value = "Vergara, Gustavo"
$ws.Range("C213").Value = "Parana"
$ws.Range("D213").Value = "Tiro Federal"

# reorder tied-rating group rows [226, 227] alphabetically by player name
$ws.Range("B226").Value = "Escalante, Samuel"
$ws.Range("C226").Value = ""
$ws.Range("D226").Value = ""
$ws.Range("B227").Value = "Musuruana, Francisco"
$ws.Range("C227").Value = ""
$ws.Range("D227").Value = ""

# reorder tied-rating group rows [241, 242] alphabetically by player name
$ws.Range("B241").Value = "Antunez, Pablo"
$ws.Range("C241").Value = ""
$ws.Range("D241").Value = ""
$ws.Range("B242").Value = "Ferrero, Alejandro"
$ws.Range("C242").Value = ""
$ws.Range("D242").Value = ""

# reorder tied-rating group rows [244, 245] alphabetically by player name
$ws.Range("B244").Value = "Lell, Claudia"
$ws.Range("C244").Value = "Parana"
$ws.Range("D244").Value = "Tiro Federal"
$ws.Range("B245").Value = "Velazquez, Pedro"
$ws.Range("C245").Value = ""
$ws.Range("D245").Value = ""

# reorder tied-rating group rows [252, 253] alphabetically by player name
$ws.Range("B252").Value = "Bertoli, Julian"
$ws.Range("C252").Value = ""
$ws.Range("D252").Value = ""
$ws.Range("B253").Value = "Bertoli, Maximiliano"
$ws.Range("C253").Value = ""
$ws.Range("D253").Value = ""

# reorder tied-rating group rows [254, 255] alphabetically by player name
$ws.Range("B254").Value = "Macor, Mateo"
$ws.Range("C254").Value = ""
$ws.Range("D254").Value = ""
$ws.Range("B255").Value = "Seib, Silvia"
$ws.Range("C255").Value = "Parana"
$ws.Range("D255").Value = "Aspatem"

